$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-removed rows (11-13), which collapsed the 4th target cluster (Resolving-Mac)
$ws.Rows("11:13").Delete()

# Row 2
$ws.Range("G2").Value = 6.156604333333333
$ws.Range("H2").Value = 18.469813
$ws.Range("I2").Value = 0.3861700262161295
$ws.Range("J2").Value = 0.3861700262161295
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.261293
$ws.Range("N2").Value = 0.783879
$ws.Range("O2").Value = 0.0361900776512412
$ws.Range("P2").Value = 0.03619007765124121
$ws.Range("Q2").Value = 1.608677616069667
$ws.Range("R2").Value = 14.478098544627
$ws.Range("S2").Value = 0.01397552323534358
$ws.Range("T2").Value = 0.01397552323534358

# Row 3
$ws.Range("G3").Value = 6.156604333333333
$ws.Range("H3").Value = 18.469813
$ws.Range("I3").Value = 0.3861700262161295
$ws.Range("J3").Value = 0.3861700262161295
$ws.Range("M3").Value = 0.7775033333333333
$ws.Range("O3").Value = 0.1076871787894517
$ws.Range("P3").Value = 0.1076871787894517
$ws.Range("Q3").Value = 4.786780391181111
$ws.Range("R3").Value = 43.08102352063
$ws.Range("S3").Value = 0.04158556065626359
$ws.Range("T3").Value = 0.0415855606562636

# Row 4
$ws.Range("G4").Value = 6.156604333333333
$ws.Range("H4").Value = 18.469813
$ws.Range("I4").Value = 0.3861700262161295
$ws.Range("J4").Value = 0.3861700262161295
$ws.Range("M4").Value = 6.181221333333333
$ws.Range("N4").Value = 18.543664
$ws.Range("O4").Value = 0.8561227435593071
$ws.Range("P4").Value = 0.8561227435593072
$ws.Range("Q4").Value = 38.05533404609244
$ws.Range("R4").Value = 342.498006414832
$ws.Range("S4").Value = 0.3306089423245223
$ws.Range("T4").Value = 0.3306089423245224

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "ECs"
$ws.Range("G5").Value = 5.867977666666667
$ws.Range("H5").Value = 17.603933
$ws.Range("I5").Value = 0.3680660582820729
$ws.Range("J5").Value = 0.3680660582820729
$ws.Range("M5").Value = 0.261293
$ws.Range("N5").Value = 0.783879
$ws.Range("O5").Value = 0.0361900776512412
$ws.Range("P5").Value = 0.03619007765124121
$ws.Range("Q5").Value = 1.533261488456334
$ws.Range("R5").Value = 13.799353396107
$ws.Range("S5").Value = 0.01332033923001449
$ws.Range("T5").Value = 0.01332033923001449

# Row 6
$ws.Range("D6").Value = "FAPs"
$ws.Range("G6").Value = 5.867977666666667
$ws.Range("I6").Value = 0.3680660582820729
$ws.Range("J6").Value = 0.3680660582820729
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7775033333333333
$ws.Range("N6").Value = 2.33251
$ws.Range("O6").Value = 0.1076871787894517
$ws.Range("P6").Value = 0.1076871787894517
$ws.Range("Q6").Value = 4.56237219575889
$ws.Range("R6").Value = 41.06134976183
$ws.Range("S6").Value = 0.03963599542455034
$ws.Range("T6").Value = 0.03963599542455035

# Row 7
$ws.Range("D7").Value = "MuSCs"
$ws.Range("G7").Value = 5.867977666666667
$ws.Range("I7").Value = 0.3680660582820729
$ws.Range("J7").Value = 0.3680660582820729
$ws.Range("M7").Value = 6.181221333333333
$ws.Range("N7").Value = 18.543664
$ws.Range("O7").Value = 0.8561227435593071
$ws.Range("P7").Value = 0.8561227435593072
$ws.Range("Q7").Value = 36.27126873672356
$ws.Range("R7").Value = 326.441418630512
$ws.Range("S7").Value = 0.3151097236275081
$ws.Range("T7").Value = 0.3151097236275082

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("D8").Value = "ECs"
$ws.Range("G8").Value = 3.918147666666667
$ws.Range("H8").Value = 11.754443
$ws.Range("I8").Value = 0.2457639155017975
$ws.Range("J8").Value = 0.2457639155017975
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.261293
$ws.Range("N8").Value = 0.783879
$ws.Range("O8").Value = 0.0361900776512412
$ws.Range("P8").Value = 0.03619007765124121
$ws.Range("Q8").Value = 1.023784558266334
$ws.Range("R8").Value = 9.214061024397001
$ws.Range("S8").Value = 0.008894215185883134
$ws.Range("T8").Value = 0.008894215185883136

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("D9").Value = "FAPs"
$ws.Range("G9").Value = 3.918147666666667
$ws.Range("H9").Value = 11.754443
$ws.Range("I9").Value = 0.2457639155017975
$ws.Range("J9").Value = 0.2457639155017975
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.7775033333333333
$ws.Range("N9").Value = 2.33251
$ws.Range("O9").Value = 0.1076871787894517
$ws.Range("P9").Value = 0.1076871787894517
$ws.Range("Q9").Value = 3.046372871325556
$ws.Range("R9").Value = 27.41735584193001
$ws.Range("S9").Value = 0.02646562270863777
$ws.Range("T9").Value = 0.02646562270863778

# Row 10
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.918147666666667
$ws.Range("H10").Value = 11.754443
$ws.Range("I10").Value = 0.2457639155017975
$ws.Range("J10").Value = 0.2457639155017975
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 6.181221333333333
$ws.Range("N10").Value = 18.543664
$ws.Range("O10").Value = 0.8561227435593071
$ws.Range("P10").Value = 0.8561227435593072
$ws.Range("Q10").Value = 24.21893794435023
$ws.Range("R10").Value = 217.970441499152
$ws.Range("S10").Value = 0.2104040776072766
$ws.Range("T10").Value = 0.2104040776072766
